# Update "想去人数" (interested-count) figures in the F column across the
# four sheets, matching the refreshed data pull baked into the diff.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1201
$ws.Range("F8").Value = 428
$ws.Range("F10").Value = 57
$ws.Range("F11").Value = 319
$ws.Range("F12").Value = 276
$ws.Range("F13").Value = 1650
$ws.Range("F16").Value = 771
$ws.Range("F19").Value = 12533
$ws.Range("F20").Value = 12580
$ws.Range("F25").Value = 43
$ws.Range("F26").Value = 448
$ws.Range("F27").Value = 1954
$ws.Range("F30").Value = 223
$ws.Range("F31").Value = 647

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 2

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 136

# 全部类型 (All types - combined view)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1201
$ws.Range("F7").Value = 136
$ws.Range("F10").Value = 428
$ws.Range("F13").Value = 57
$ws.Range("F15").Value = 319
$ws.Range("F17").Value = 276
$ws.Range("F18").Value = 1650
$ws.Range("F21").Value = 771
$ws.Range("F25").Value = 12533
$ws.Range("F26").Value = 12580
$ws.Range("F31").Value = 43
$ws.Range("F32").Value = 448
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = 1954
$ws.Range("F40").Value = 223
$ws.Range("F41").Value = 647
